$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header for column C: audioFalse -> currentPhase
$ws.Range("C1").Value = "currentPhase"

# Collapse the per-row audio file references in column C into a constant "train1P2"
$ws.Range("C2").Value = "train1P2"
$ws.Range("C3").Value = "train1P2"
